# Login screen done I think
# Add a new inventory item ("Talco para bebé") to the "Inventario" sheet,
# and normalize the pre-existing rows 5-7 (previously stored with a
# floating-point literal such as 36.0) back to plain whole numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventario")
$ws.Select()

# Normalize the pre-existing price/stock figures for rows 5-7.
$ws.Range("C5").Value = 36
$ws.Range("D5").Value = 30
$ws.Range("C6").Value = 30
$ws.Range("D6").Value = 50
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 5

# New row 8: Folio 3881, "Talco para bebé", Precio 30, Stock 20.
# The leading apostrophe forces the numeric-looking Folio to be stored
# as text (matching the existing Folio column, which is text everywhere
# else), then the style is reset back to Normal so no extra quote-prefix
# cell style is left behind.
$ws.Range("A8").Value = "'3881"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = "Talco para bebé"
$ws.Range("C8").Value = 30
$ws.Range("D8").Value = 20

# Match the new selection left behind on the sheet.
$ws.Range("A9:D10").Select()
